$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R: year header 2021 (copy format from Q4, same header style)
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# New column R: data value for 2021 (copy format from H5, same 0.0 numeric style)
$ws.Range("H5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 18.953297329007047

# Move the active selection from Q9 to Q8
$ws.Range("Q8").Select()
